# Applies the "started to implement ressource grid subscriber and axi ring
# writer" edit: re-positions a handful of connectors/labels in the small
# diagram cluster around (9.3M, 4.0-4.7M EMU), renames the "busy" label to
# "status" (and widens its box), and removes the now-redundant "bram" and
# "ring#" labels together with the stray arrow connector that pointed at
# "ring#".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU <-> point conversion on this COM host is lossy (Left/Top/Width/Height
# are single-precision floats under the hood, and the EMU value is derived
# by truncating, not rounding). Nudge the point value up in tiny steps
# until the round-tripped EMU matches the desired integer exactly.
$EMU_PER_PT = 12700.0
$STEP_PT = 0.00001
$MAX_ITER = 1000

function Set-PreciseLeft($shape, [double]$targetEmu) {
    $ptsVal = $targetEmu / $EMU_PER_PT
    for ($i = 0; $i -lt $MAX_ITER; $i++) {
        $shape.Left = $ptsVal
        $gotEmu = [Math]::Round($shape.Left * $EMU_PER_PT)
        if ($gotEmu -eq $targetEmu) {
            break
        }
        $ptsVal = $ptsVal + $STEP_PT
    }
}

function Set-PreciseTop($shape, [double]$targetEmu) {
    $ptsVal = $targetEmu / $EMU_PER_PT
    for ($i = 0; $i -lt $MAX_ITER; $i++) {
        $shape.Top = $ptsVal
        $gotEmu = [Math]::Round($shape.Top * $EMU_PER_PT)
        if ($gotEmu -eq $targetEmu) {
            break
        }
        $ptsVal = $ptsVal + $STEP_PT
    }
}

function Set-PreciseWidth($shape, [double]$targetEmu) {
    $ptsVal = $targetEmu / $EMU_PER_PT
    for ($i = 0; $i -lt $MAX_ITER; $i++) {
        $shape.Width = $ptsVal
        $gotEmu = [Math]::Round($shape.Width * $EMU_PER_PT)
        if ($gotEmu -eq $targetEmu) {
            break
        }
        $ptsVal = $ptsVal + $STEP_PT
    }
}

function Set-PreciseHeight($shape, [double]$targetEmu) {
    $ptsVal = $targetEmu / $EMU_PER_PT
    for ($i = 0; $i -lt $MAX_ITER; $i++) {
        $shape.Height = $ptsVal
        $gotEmu = [Math]::Round($shape.Height * $EMU_PER_PT)
        if ($gotEmu -eq $targetEmu) {
            break
        }
        $ptsVal = $ptsVal + $STEP_PT
    }
}

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cand = $shapes.Item($i)
        if ($cand.Id -eq $id) {
            return $cand
        }
    }
    return $null
}

# 1) "Straight Arrow Connector 87" (id 88): slide its y offset down.
$shp88 = Get-ShapeById $s.Shapes 88
Set-PreciseTop $shp88 4087593

# 2) "Straight Arrow Connector 20" (id 77, flipH): move both x/y.
$shp77 = Get-ShapeById $s.Shapes 77
Set-PreciseLeft $shp77 9332123
Set-PreciseTop $shp77 4351201

# 3) "TextBox 61" (id 90, text "busy" -> "status"): reposition, widen, retext.
$shp90 = Get-ShapeById $s.Shapes 90
Set-PreciseLeft $shp90 9279523
Set-PreciseTop $shp90 4095752
Set-PreciseWidth $shp90 557973
Set-PreciseHeight $shp90 276999
$shp90.TextFrame.TextRange.Text = "status"

# 4) "TextBox 61" (id 104, text "bram"): removed entirely.
$shp104 = Get-ShapeById $s.Shapes 104
$shp104.Delete()

# 5) "Straight Arrow Connector 20" (id 113): move both x/y.
$shp113 = Get-ShapeById $s.Shapes 113
Set-PreciseLeft $shp113 9344292
Set-PreciseTop $shp113 4610132

# 6) "TextBox 61" (id 122, text "start"): reposition only (size/text unchanged).
$shp122 = Get-ShapeById $s.Shapes 122
Set-PreciseLeft $shp122 9290064
Set-PreciseTop $shp122 4375123

# 7) "Straight Arrow Connector 20" (id 126, flipH): removed entirely.
$shp126 = Get-ShapeById $s.Shapes 126
$shp126.Delete()

# 8) "TextBox 61" (id 128, text "ring#"): removed entirely.
$shp128 = Get-ShapeById $s.Shapes 128
$shp128.Delete()
